$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-14 08:18:22"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-14 08:18:11"
$wsZhCn.Range("K2").Value = "2016-10-14 08:18:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-10-14 08:18:22"
$wsDeDe.Range("K2").Value = "2016-10-14 08:19:07"
